$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "CamOffestPos"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "string"
$ws.Range("C13").Value = $false
$ws.Range("D13").Value = $false
$ws.Range("E13").Value = $false
$ws.Range("F13").Value = $true
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "Friend"
$ws.Range("J13").NumberFormat = "@"
$ws.Range("J13").Value = "acctorid"

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "CamOffestRot"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "string"
$ws.Range("C14").Value = $false
$ws.Range("D14").Value = $false
$ws.Range("E14").Value = $false
$ws.Range("F14").Value = $true
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "Friend"
$ws.Range("J14").NumberFormat = "@"
$ws.Range("J14").Value = "acctorid"

$ws.Range("A14").Select()
